$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(75, 1).Value = 44388.7674169618
$ws.Cells.Item(76, 1).Value = 44389.76861943601
$ws.Cells.Item(76, 2).Value = 79668
$ws.Cells.Item(76, 3).Value = 67217
$ws.Cells.Item(76, 4).Value = 3575
$ws.Cells.Item(76, 5).Value = 2205
$ws.Cells.Item(76, 6).Value = 1579
$ws.Cells.Item(76, 7).Value = 21083
$ws.Cells.Item(76, 8).Value = 1556
$ws.Cells.Item(76, 9).Value = 899
$ws.Cells.Item(76, 10).Value = 197
